# Applies the CTDC pubmed id / trial arm update:
#  - Adds a new "TabName" column (A) identifying each query row as CasesTab / FilesTab
#  - Replaces the old single Cases query with updated Cases + new Files queries
#  - Replaces the old stat query with an updated stat query, duplicated for both tabs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column; existing columns A:D shift to B:E
$ws.Columns.Item(1).Insert()

# ---- Row 1 : headers ----
$ws.Range("A1").Value = 'TabName'
$ws.Range("B1").Value = 'query'
$ws.Range("C1").Value = 'StatQuery'
$ws.Range("D1").Value = 'dbExcel'
$ws.Range("E1").Value = 'WebExcel'

# ---- Row 2 : Cases tab ----
$ws.Range("A2").Value = 'CasesTab'
$ws.Range("B2").Value = 'MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
    WHERE a.pubmed_id IN [''31765263''] 
OPTIONAL MATCH (f:file)-[*]->(c)
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity'
$ws.Range("C2").Value = 'MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
  WHERE a.pubmed_id IN [''31765263'']
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files'
$ws.Range("D2").Value = 'TC02_Trials_Filter_PubmedID-317_Neo4jData.xlsx'
$ws.Range("E2").Value = 'TC02_Trials_Filter_PubmedID-317_WebData.xlsx'
$ws.Range("B2:C2").WrapText = $true

# ---- Row 3 : Files tab (new row) ----
$ws.Range("A3").Value = 'FilesTab'
$ws.Range("B3").Value = 'MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
WHERE a.pubmed_id IN [''31765263'']
WITH
    f, parent, c, a, ct,
    [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`'
$ws.Range("C3").Value = 'MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
  WHERE a.pubmed_id IN [''31765263'']
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files'
$ws.Range("D3").Value = 'TC02_Trials_Filter_PubmedID-317_Neo4jData.xlsx'
$ws.Range("E3").Value = 'TC02_Trials_Filter_PubmedID-317_WebData.xlsx'
$ws.Range("B3:C3").WrapText = $true

# ---- Row heights for the wrapped, multi-line query cells ----
$ws.Rows.Item(2).RowHeight = 195
$ws.Rows.Item(3).RowHeight = 409.5

# ---- Column widths ----
$ws.Columns.Item(1).ColumnWidth = 8
$ws.Columns.Item(2).ColumnWidth = 75
$ws.Columns.Item(3).ColumnWidth = 75
$ws.Columns.Item(4).ColumnWidth = 69.5
$ws.Columns.Item(5).ColumnWidth = 27.666666666666668

# ---- Active selection ----
$ws.Range("C3").Select()
